# Add PF/1.0.5 to meta-sheet
# New row 3: version "PF/1.0.5" in column A, with "X" markers for the
# dev2 / sit2 / uat2 / prod columns (B, C, D).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "PF/1.0.5"
$ws.Range("B3").Value = "X"
$ws.Range("C3").Value = "X"
$ws.Range("D3").Value = "X"

# Keep the new row on the sheet's default (unstyled) format rather than
# inheriting the bold/aligned style used by the header & existing data rows.
$ws.Range("A3:D3").Style = "Normal"
